# Apply manuscript review updates to Table 2 (age statistics for sampled moraines)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Round the reported ages/uncertainties to one decimal place ---

# Tallada
$ws.Range("G4").Value = "3.2 ± 0.7"

# Outer Pleta Naua (keep trailing "h" footnote marker superscript)
$ws.Range("G5").Value = "12.5 ± 0.4h"
$ws.Range("G5").Characters(11,1).Font.Superscript = $true

# Aranser (keep trailing "i" footnote marker superscript)
$ws.Range("G6").Value = "23.3 ± 1.1i"
$ws.Range("G6").Characters(11,1).Font.Superscript = $true

# Soum d'Ech - Right
$ws.Range("G7").Value = "22.3 ± 0.9"

# Soum d'Ech - Outer
$ws.Range("G8").Value = "26.2 ± 2.5"

# Soum d'Ech - Inner
$ws.Range("G9").Value = "26.1 ± 1.7"

# Soum d'Ech - Combined
$ws.Range("G10").Value = "27.3 ± 1.8"

# --- Update the view: zoom to 130% and select the footnote row ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 130
$ws.Range("B11:M11").Select()
